$wb = $excel.ActiveWorkbook

# compraProducto: B2 "Si" -> "No", selection moves to B2
$ws = $wb.Worksheets.Item("compraProducto")
$ws.Range("B2").Value = "No"
$ws.Range("B2").Select() | Out-Null

# crearCuenta: B2 "Si" -> "No", selection moves to B2
$ws = $wb.Worksheets.Item("crearCuenta")
$ws.Range("B2").Value = "No"
$ws.Range("B2").Select() | Out-Null

# alertsTools: B2 "Si" -> "No", selection moves to B2
$ws = $wb.Worksheets.Item("alertsTools")
$ws.Range("B2").Value = "No"
$ws.Range("B2").Select() | Out-Null

# dateTools: no value change, only selection moves to B2
$ws = $wb.Worksheets.Item("dateTools")
$ws.Range("B2").Select() | Out-Null

# mercuryTours (active sheet): B2 "Si" -> "No", selection moves to C9
$ws = $wb.Worksheets.Item("mercuryTours")
$ws.Range("B2").Value = "No"
$ws.Range("C9").Select() | Out-Null

# opciones: no value change, only selection moves to H17
$ws = $wb.Worksheets.Item("opciones")
$ws.Range("H17").Select() | Out-Null

# Keep mercuryTours as the active/visible sheet with C9 selected
$wsActive = $wb.Worksheets.Item("mercuryTours")
$wsActive.Activate()
$wsActive.Range("C9").Select() | Out-Null
